$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text ("@") before assigning so that Excel
# preserves the exact string representation (e.g. thousand-separator-
# style prices like "69.144.49", leading/trailing zeros, percentages)
# instead of auto-converting numeric-looking strings into numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.144.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.839.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.828.68"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.85%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.476"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.02"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000249"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.481.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.844.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.371.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.41"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.31"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "498.92"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.54"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.55%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.52"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000138"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.81%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.49"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.03"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +9.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.82"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.01"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "461.92"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.330"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.04"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.43"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.887.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.38"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.64"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.13%  "
